$wb = $excel.ActiveWorkbook

# --- Sheet: Population Definitions ---
$popDef = $wb.Worksheets.Item("Population Definitions")
$popDef.Range("A2").Value = "Children"
$popDef.Range("A3").Value = "General Population"
$popDef.Range("B2").Value = "SAC"
$popDef.Range("B3").Value = "GEN"
$popDef.Activate()
$popDef.Range("B4").Select()

# --- Sheet: Population Sizes ("Number of births" block, rows 5-7) ---
$popSizes = $wb.Worksheets.Item("Population Sizes")
$popSizes.Range("F6").Value = 500
$popSizes.Range("G6").Value = 600
$popSizes.Range("G6").Style = "Normal"
$popSizes.Range("H6").Value = 700
$popSizes.Range("I6").Value = 800
$popSizes.Range("J6").Value = 900
$popSizes.Range("J6").Style = "Normal"

$popSizes.Range("C7").Value = 0
$popSizes.Range("E7").Value = $null
$popSizes.Range("F7").Value = $null
$popSizes.Range("G7").Value = $null
$popSizes.Activate()
$popSizes.Range("J6").Select()

# --- Sheet: Epidemic Characteristics ---
$epi = $wb.Worksheets.Item("Epidemic Characteristics")
$epi.Activate()
$epi.Range("C3").Select()

# --- Sheet: Cascade Parameters (becomes the active sheet / tab) ---
$cascade = $wb.Worksheets.Item("Cascade Parameters")
$cascade.Activate()
